$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# "Preschool: prep desample for school, integration of alt forms"
# The Child row (age_range=Child, form=Home / School) and the Teen row
# (age_range=Teen, form=Home / School / Self) are being folded into the
# norms-output the same way the Adult rows already were: flip their
# existing "X" markers (data exists, unused) to "XX" (data exists AND is
# integrated into norms-output) for every source column that has a mark.
$cellsToMark = @(
    "C4", "D4", "E4", "G4", "H4",   # row 4  - Child / Home
    "C5", "G5",                      # row 5  - Child / School
    "C11", "D11", "E11", "G11", "H11", # row 11 - Teen / Home
    "C12", "G12",                    # row 12 - Teen / School
    "C13", "D13", "E13", "G13", "H13"  # row 13 - Teen / Self
)

foreach ($cellRef in $cellsToMark) {
    $ws.Range($cellRef).Value = "XX"
}

# Reflect the author's last-clicked cell in the saved view.
$ws.Range("H14").Select()
